$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new data row (row 2) with the weight-observation values.
$ws.Range("A2").Value = "weight-observation"
$ws.Range("B2").Value = "Weight Observation"
$ws.Range("E2").Value = "LOINC#29463-7"
$ws.Range("G2").Value = "dateTime, Period, Timing, instant"
$ws.Range("H2").Value = "Quantity" + [char]0x0135
$ws.Range("I2").Value = "optional"

# Copy the header row's formatting onto the new row so it keeps the same style.
$ws.Range("A1:K1").Copy()
$ws.Range("A2:K2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
